$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1264559.9
$ws.Range("I31").Value = 1264559.9
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3793679.7
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3793449.7
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 97.625
$ws.Range("I39").Value = 54.42857
$ws.Range("J39").Value = 400
$ws.Range("K39").Value = 163.28571
$ws.Range("L39").Value = 1200
$ws.Range("M39").Value = 132.71429
$ws.Range("N39").Value = -1792

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2785.077
$ws.Range("I62").Value = 2242.8572
$ws.Range("J62").Value = 3417.6667
$ws.Range("K62").Value = 2242.8572
$ws.Range("L62").Value = 3417.6667
$ws.Range("M62").Value = -1618.8572
$ws.Range("N62").Value = -4665.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2785.077
$ws.Range("I65").Value = 2242.8572
$ws.Range("J65").Value = 3417.6667
$ws.Range("K65").Value = 11214.286
$ws.Range("L65").Value = 17088.3335
$ws.Range("M65").Value = -8094.286
$ws.Range("N65").Value = -23328.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3612.875
$ws.Range("I76").Value = 3500.5
$ws.Range("J76").Value = 3950
$ws.Range("K76").Value = 3500.5
$ws.Range("L76").Value = 3950
$ws.Range("M76").Value = -3185.5
$ws.Range("N76").Value = -4580

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3612.875
$ws.Range("I79").Value = 3500.5
$ws.Range("J79").Value = 3950
$ws.Range("K79").Value = 3500.5
$ws.Range("L79").Value = 3950
$ws.Range("M79").Value = -2408.5
$ws.Range("N79").Value = -6134

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 41360.6
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 41360.6
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 41360.6
$ws.Range("N93").Value = -46352.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1009.5179
$ws.Range("I129").Value = 596.6667
$ws.Range("J129").Value = 1032.8868
$ws.Range("K129").Value = 1790.0001
$ws.Range("L129").Value = 3098.6604
$ws.Range("M129").Value = 3209.9999
$ws.Range("N129").Value = -13098.6604

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2872.186
$ws.Range("I138").Value = 1378.5
$ws.Range("J138").Value = 3324.818
$ws.Range("K138").Value = 4135.5
$ws.Range("L138").Value = 9974.454000000002
$ws.Range("M138").Value = 1004.5
$ws.Range("N138").Value = -20254.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9466.029
$ws.Range("I32").Value = 9764.462
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 9764.462
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -9477.462
$ws.Range("N32").Value = -3574

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1674.3
$ws.Range("I61").Value = 1551.8948
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1551.8948
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1339.8948
$ws.Range("N61").Value = -4424

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 70933
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 70933
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 70933
$ws.Range("N92").Value = -75925

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2242.8572
$ws.Range("I122").Value = 1989.2
$ws.Range("J122").Value = 2383.7778
$ws.Range("K122").Value = 5967.6
$ws.Range("L122").Value = 7151.3334
$ws.Range("M122").Value = -3517.6
$ws.Range("N122").Value = -12051.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2036.75
$ws.Range("I132").Value = 1410.5
$ws.Range("J132").Value = 3498
$ws.Range("K132").Value = 4231.5
$ws.Range("L132").Value = 10494
$ws.Range("M132").Value = -1701.5
$ws.Range("N132").Value = -15554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1674.3
$ws.Range("I136").Value = 1551.8948
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 4655.6844
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -2105.6844
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 8599.6
$ws.Range("I5").Value = 5999.3335
$ws.Range("J5").Value = 12500
$ws.Range("K5").Value = 5999.3335
$ws.Range("L5").Value = 12500
$ws.Range("M5").Value = -5886.3335
$ws.Range("N5").Value = -12726

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5367
$ws.Range("I36").Value = 3524
$ws.Range("J36").Value = 9053
$ws.Range("K36").Value = 3524
$ws.Range("L36").Value = 9053
$ws.Range("M36").Value = -3136
$ws.Range("N36").Value = -9829

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 15000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 15000
$ws.Range("N37").Value = -15214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 5367
$ws.Range("I40").Value = 3524
$ws.Range("J40").Value = 9053
$ws.Range("K40").Value = 3524
$ws.Range("L40").Value = 9053
$ws.Range("M40").Value = -3364
$ws.Range("N40").Value = -9373

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 30000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 30000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30856
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 12000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 12000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12952

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1456.1
$ws.Range("I58").Value = 1678.2858
$ws.Range("J58").Value = 937.6667
$ws.Range("K58").Value = 1678.2858
$ws.Range("L58").Value = 937.6667
$ws.Range("M58").Value = -1475.2858
$ws.Range("N58").Value = -1343.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3092.75
$ws.Range("I122").Value = 3262.4614
$ws.Range("J122").Value = 2892.182
$ws.Range("K122").Value = 9787.3842
$ws.Range("L122").Value = 8676.545999999998
$ws.Range("M122").Value = -7337.3842
$ws.Range("N122").Value = -13576.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1778.2632
$ws.Range("I132").Value = 1391.5385
$ws.Range("J132").Value = 2616.1667
$ws.Range("K132").Value = 4174.6155
$ws.Range("L132").Value = 7848.500100000001
$ws.Range("M132").Value = -1644.6155
$ws.Range("N132").Value = -12908.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1456.1
$ws.Range("I136").Value = 1678.2858
$ws.Range("J136").Value = 937.6667
$ws.Range("K136").Value = 5034.857400000001
$ws.Range("L136").Value = 2813.0001
$ws.Range("M136").Value = -2484.857400000001
$ws.Range("N136").Value = -7913.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 993
$ws.Range("I42").Value = 993
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 2979
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -2445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1800
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1800
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5400
$ws.Range("N51").Value = -6320
$ws.Range("M51").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 3850
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3850
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 11550
$ws.Range("N58").Value = -11806

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 742.1539
$ws.Range("I92").Value = 799.6667
$ws.Range("J92").Value = 692.8570999999999
$ws.Range("K92").Value = 2399.0001
$ws.Range("L92").Value = 2078.5713
$ws.Range("M92").Value = -1151.0001
$ws.Range("N92").Value = -4574.5713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 7844.125
$ws.Range("I109").Value = 10732.6
$ws.Range("J109").Value = 3030
$ws.Range("K109").Value = 32197.8
$ws.Range("L109").Value = 9090
$ws.Range("M109").Value = -31157.8
$ws.Range("N109").Value = -11170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 2099.2
$ws.Range("I115").Value = 548.6667
$ws.Range("J115").Value = 2763.7144
$ws.Range("K115").Value = 1646.0001
$ws.Range("L115").Value = 8291.143199999999
$ws.Range("M115").Value = -471.0001
$ws.Range("N115").Value = -10641.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 604.3333
$ws.Range("I122").Value = 400
$ws.Range("J122").Value = 672.44446
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 6052.00014
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -10952.00014

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5326.923
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 5354.1665
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 16062.4995
$ws.Range("M133").Value = -9940
$ws.Range("N133").Value = -26182.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4519.341
$ws.Range("I134").Value = 1790
$ws.Range("J134").Value = 5931.069
$ws.Range("K134").Value = 5370
$ws.Range("L134").Value = 17793.207
$ws.Range("M134").Value = -300
$ws.Range("N134").Value = -27933.207

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2071.5
$ws.Range("I138").Value = 896.6667
$ws.Range("J138").Value = 2575
$ws.Range("K138").Value = 2690.0001
$ws.Range("L138").Value = 7725
$ws.Range("M138").Value = 2449.9999
$ws.Range("N138").Value = -18005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1782.64
$ws.Range("I139").Value = 1187.5
$ws.Range("J139").Value = 2840.6667
$ws.Range("K139").Value = 3562.5
$ws.Range("L139").Value = 8522.000100000001
$ws.Range("M139").Value = 1577.5
$ws.Range("N139").Value = -18802.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2711.7334
$ws.Range("I140").Value = 1181
$ws.Range("J140").Value = 5773.2
$ws.Range("K140").Value = 3543
$ws.Range("L140").Value = 17319.6
$ws.Range("M140").Value = 1637
$ws.Range("N140").Value = -27679.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1000.6667
$ws.Range("I41").Value = 1000.6667
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1000.6667
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -645.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1643.3529
$ws.Range("I122").Value = 1295.1538
$ws.Range("J122").Value = 2775
$ws.Range("K122").Value = 3885.4614
$ws.Range("L122").Value = 8325
$ws.Range("M122").Value = -1435.4614
$ws.Range("N122").Value = -13225

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 650.5
$ws.Range("I55").Value = 301
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 301
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -128
$ws.Range("N55").Value = -1346

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 28750
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 28750
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 28750
$ws.Range("N92").Value = -33742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10420948
$ws.Range("I122").Value = 14709597
$ws.Range("J122").Value = 5657.143
$ws.Range("K122").Value = 44128791
$ws.Range("L122").Value = 16971.429
$ws.Range("M122").Value = -44126341
$ws.Range("N122").Value = -21871.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 25150
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 49800
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 49800
$ws.Range("M8").Value = -360
$ws.Range("N8").Value = -50080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 32500
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 32500
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 32500
$ws.Range("N105").Value = -39488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6759130.5
$ws.Range("I122").Value = 11365586
$ws.Range("J122").Value = 2995.2
$ws.Range("K122").Value = 34096758
$ws.Range("L122").Value = 8985.599999999999
$ws.Range("M122").Value = -34094308
$ws.Range("N122").Value = -13885.6
